{"js": "// Apply Word formatting to style by default: give each \"column\" of the\n// verse tables (and the copyright table) its own paragraph style, instead\n// of sharing a single MSCParagraph / no-style / MSCCopyright style across\n// columns. New custom styles MSC_Paragraph_A/B/C are created (based on the\n// existing MSCParagraph style) and assigned to the relevant paragraphs.\n\n// 1) Create the three new paragraph styles, each based on MSCParagraph.\nconst newStyleNames = [\"MSC_Paragraph_A\", \"MSC_Paragraph_B\", \"MSC_Paragraph_C\"];\nfor (const name of newStyleNames) {\n  context.document.addStyle(name, Word.StyleType.paragraph);\n}\nawait context.sync();\n\n// addStyle()'s returned anchor does not reliably round-trip property\n// writes, so re-look the styles up by name before configuring them.\nconst styles = context.document.getStyles();\nconst styleA = styles.getByNameOrNullObject(\"MSC_Paragraph_A\");\nconst styleB = styles.getByNameOrNullObject(\"MSC_Paragraph_B\");\nconst styleC = styles.getByNameOrNullObject(\"MSC_Paragraph_C\");\nstyleA.baseStyle = \"MSCParagraph\";\nstyleB.baseStyle = \"MSCParagraph\";\nstyleC.baseStyle = \"MSCParagraph\";\nawait context.sync();\n\n// 2) Walk the two verse tables (Genesis, Mark). In each, the body row is\n// the last row (Genesis has a header \"version name\" row first; Mark does\n// not). Column A (cell 0) keeps the existing MSCParagraph-styled\n// paragraphs but retargets them to MSC_Paragraph_A; columns B and C\n// (cells 1 and 2) are plain/unstyled paragraphs that get MSC_Paragraph_B\n// / MSC_Paragraph_C respectively.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const t of tables.items) {\n  t.load(\"rowCount\");\n}\nawait context.sync();\n\nconst verseTables = tables.items.slice(0, 2);\n\nfor (const t of verseTables) {\n  const bodyRowIndex = t.rowCount - 1;\n\n  const cellA = t.getCell(bodyRowIndex, 0);\n  const cellB = t.getCell(bodyRowIndex, 1);\n  const cellC = t.getCell(bodyRowIndex, 2);\n\n  const parasA = cellA.body.paragraphs;\n  const parasB = cellB.body.paragraphs;\n  const parasC = cellC.body.paragraphs;\n  parasA.load(\"items\");\n  parasB.load(\"items\");\n  parasC.load(\"items\");\n  await context.sync();\n\n  for (const p of parasA.items) {\n    p.load(\"style\");\n  }\n  await context.sync();\n\n  for (const p of parasA.items) {\n    if (p.style === \"MSC_Paragraph\") {\n      p.style = \"MSC_Paragraph_A\";\n    }\n  }\n  for (const p of parasB.items) {\n    p.style = \"MSC_Paragraph_B\";\n  }\n  for (const p of parasC.items) {\n    p.style = \"MSC_Paragraph_C\";\n  }\n  await context.sync();\n}\n\n// 3) Copyright table: the license text cell already uses MSCCopyright;\n// give the other two (empty) cells in that row the same style.\nconst copyrightTable = tables.items[2];\ncopyrightTable.load(\"rowCount\");\nawait context.sync();\n\nconst lastCopyrightRow = copyrightTable.rowCount - 1;\nconst copyB = copyrightTable.getCell(lastCopyrightRow, 1);\nconst copyC = copyrightTable.getCell(lastCopyrightRow, 2);\nconst copyParasB = copyB.body.paragraphs;\nconst copyParasC = copyC.body.paragraphs;\ncopyParasB.load(\"items\");\ncopyParasC.load(\"items\");\nawait context.sync();\n\nfor (const p of copyParasB.items) {\n  p.style = \"MSCCopyright\";\n}\nfor (const p of copyParasC.items) {\n  p.style = \"MSCCopyright\";\n}\nawait context.sync();\n", "ps1": "# Apply Word formatting to style by default: give each \"column\" of the\n# verse tables (and the copyright table) its own paragraph style, instead\n# of sharing a single MSCParagraph / no-style / MSCCopyright style across\n# columns. New custom styles MSC_Paragraph_A/B/C are created (based on the\n# existing MSCParagraph style) and assigned to the relevant paragraphs.\n\n$d = $word.ActiveDocument\n\n# 1) Create the three new paragraph styles, each based on MSCParagraph.\n$styleA = $d.Styles.Add(\"MSC_Paragraph_A\", 1)\n$styleA.BaseStyle = $d.Styles(\"MSCParagraph\")\n$styleB = $d.Styles.Add(\"MSC_Paragraph_B\", 1)\n$styleB.BaseStyle = $d.Styles(\"MSCParagraph\")\n$styleC = $d.Styles.Add(\"MSC_Paragraph_C\", 1)\n$styleC.BaseStyle = $d.Styles(\"MSCParagraph\")\n\n# 2) Walk the two verse tables (Genesis, Mark). In each, the body row is\n# the last row (Genesis has a header \"version name\" row first; Mark does\n# not). Column A (cell 1) keeps the existing MSCParagraph-styled\n# paragraphs but retargets them to MSC_Paragraph_A; columns B and C\n# (cells 2 and 3) are plain/unstyled paragraphs that get MSC_Paragraph_B\n# / MSC_Paragraph_C respectively.\nfor ($ti = 1; $ti -le 2; $ti++) {\n    $t = $d.Tables($ti)\n    $bodyRow = $t.Rows.Count\n\n    $cellA = $t.Cell($bodyRow, 1)\n    $countA = $cellA.Range.Paragraphs.Count\n    for ($i = 1; $i -le $countA; $i++) {\n        $p = $cellA.Range.Paragraphs($i)\n        if ($p.Style.NameLocal -eq \"MSC_Paragraph\") {\n            $p.Style = $d.Styles(\"MSC_Paragraph_A\")\n        }\n    }\n\n    $cellB = $t.Cell($bodyRow, 2)\n    $countB = $cellB.Range.Paragraphs.Count\n    for ($i = 1; $i -le $countB; $i++) {\n        $p = $cellB.Range.Paragraphs($i)\n        $p.Style = $d.Styles(\"MSC_Paragraph_B\")\n    }\n\n    $cellC = $t.Cell($bodyRow, 3)\n    $countC = $cellC.Range.Paragraphs.Count\n    for ($i = 1; $i -le $countC; $i++) {\n        $p = $cellC.Range.Paragraphs($i)\n        $p.Style = $d.Styles(\"MSC_Paragraph_C\")\n    }\n}\n\n# 3) Copyright table: the license text cell already uses MSCCopyright;\n# give the other two (empty) cells in that row the same style.\n$copyTable = $d.Tables(3)\n$copyRow = $copyTable.Rows.Count\n\n$copyB = $copyTable.Cell($copyRow, 2)\n$countCopyB = $copyB.Range.Paragraphs.Count\nfor ($i = 1; $i -le $countCopyB; $i++) {\n    $p = $copyB.Range.Paragraphs($i)\n    $p.Style = $d.Styles(\"MSCCopyright\")\n}\n\n$copyC = $copyTable.Cell($copyRow, 3)\n$countCopyC = $copyC.Range.Paragraphs.Count\nfor ($i = 1; $i -le $countCopyC; $i++) {\n    $p = $copyC.Range.Paragraphs($i)\n    $p.Style = $d.Styles(\"MSCCopyright\")\n}\n"}
